$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fill in the run of empty paragraphs after "...AL AREA USUARIA FINAL."
#    with the new "OBJETIVO DE NEGOCIO" / "ACTORES DE NEGOCIO" blocks, while
#    keeping the surrounding blank paragraphs blank (per the diff).
# ---------------------------------------------------------------------------

# Locate the first of the run of blank paragraphs: it immediately follows
# the paragraph ending in "SE INFORMARA AL AREA USUARIA FINAL."
$anchorRange = $d.Content
$anchorFound = $anchorRange.Find.Execute("SE INFORMARA AL AREA USUARIA FINAL.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $anchorRange.Paragraphs.Item(1).Index
$firstBlankIndex = $anchorIndex + 1

# Paragraph #1 of the blank run stays blank -- insert the 3 new paragraphs
# right after it.
$p = $d.Paragraphs.Item($firstBlankIndex)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBlankIndex + 1).Range.Text = "--OBJETIVO DE NEGOCIO"

$p = $d.Paragraphs.Item($firstBlankIndex + 1)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBlankIndex + 2).Range.Text = "Disminuir el tiempo de respuesta de cada solicitud en un 50% al año anterior"

$p = $d.Paragraphs.Item($firstBlankIndex + 2)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBlankIndex + 3).Range.Text = "Incrementar  en un 30% la eficacia de verificación de lo solicitado en el stock"

# Paragraph #2 of the original blank run (now shifted by +3) stays blank.
# Paragraphs #3,#4,#5 of the original blank run become the ACTORES block.
$actStart = $firstBlankIndex + 3 + 1
$d.Paragraphs.Item($actStart).Range.Text = "--ACTORES DE NEGOCIO"
$d.Paragraphs.Item($actStart + 1).Range.Text = "Subgerencia de almacén"
$d.Paragraphs.Item($actStart + 2).Range.Text = "Área usuaria final "

# Paragraph #6 of the original blank run (the extra one) is removed so the
# final blank paragraph (#7 of the original run) is the only one left before
# "ENTREGA DE BIENES AL USUARIO".
$d.Paragraphs.Item($actStart + 3).Range.Delete()

# ---------------------------------------------------------------------------
# 2) Record a lastRenderedPageBreak right before the run of text
#    "Asistente administrativo" (last w:rPr stays untouched).
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("Asistente administrativo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $rng.Paragraphs.Item(1)
$pXml = $target.Range.WordOpenXML
$pXml = $pXml.Replace("<w:t>Asistente administrativo</w:t>", "<w:lastRenderedPageBreak/><w:t>Asistente administrativo</w:t>")

$insPoint = $target.Range.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertXML($pXml)
